# ---------------------------------------------------------------------------
# "revisions: Ch. 1 & Ch. 2"
#
# The writing-progress log (sheet "writing", table "Table1") gets five new
# daily rows (38-42, covering 2020-12-09 .. 2020-12-13). The dashboard
# chart's series ranges grow to match, the chart is resized a bit larger,
# and the selections on both sheets reflect where the author last clicked.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsWriting   = $wb.Worksheets.Item("writing")
$wsDashboard = $wb.Worksheets.Item("dashboard")

# New rows of data (column order matches the sheet: A..N)
#   A=Date B=CV C=Abstract D=Acknowledgments E=Abbreviations F=Languages
#   G=Conventions H=1:Intro I=2:Background J=3:Data&Methods K=4:Results
#   L=5:Conclusion  M=Total (=SUM(C:L))  N=Daily (=SUM(ABS(delta vs prior row)))
$newRows = @(
  @{ Row=38; Date=44174; B=568; C=239; D=87; E=198; F=492; G=416; H=7945;  I=19649; J=120; K=116; L=186 },
  @{ Row=39; Date=44175; B=568; C=239; D=87; E=198; F=492; G=416; H=8074;  I=19649; J=120; K=116; L=186 },
  @{ Row=40; Date=44176; B=568; C=239; D=87; E=198; F=492; G=416; H=8074;  I=19786; J=120; K=116; L=186 },
  @{ Row=41; Date=44177; B=568; C=239; D=87; E=198; F=492; G=416; H=8506;  I=19865; J=120; K=116; L=186 },
  @{ Row=42; Date=44178; B=568; C=239; D=87; E=198; F=492; G=416; H=8506;  I=19864; J=120; K=116; L=186 }
)

# Write all of the new cell values & formulas FIRST (while the table is
# still A1:N37) and only resize the ListObject afterwards -- doing it in the
# other order confuses array-formula assignment on the freshly-added rows.
foreach ($nr in $newRows) {
  $r = $nr.Row
  $prevRow = $r - 1

  # Column A: copy the cell above so the date number format (m/d/yyyy)
  # carries over, then overwrite with the real value.
  $wsWriting.Range("A$prevRow").Copy($wsWriting.Range("A$r"))
  $wsWriting.Range("A$r").Value2 = $nr.Date

  $wsWriting.Range("B$r").Value2 = $nr.B
  $wsWriting.Range("C$r").Value2 = $nr.C
  $wsWriting.Range("D$r").Value2 = $nr.D
  $wsWriting.Range("E$r").Value2 = $nr.E
  $wsWriting.Range("F$r").Value2 = $nr.F
  $wsWriting.Range("G$r").Value2 = $nr.G
  $wsWriting.Range("H$r").Value2 = $nr.H
  $wsWriting.Range("I$r").Value2 = $nr.I
  $wsWriting.Range("J$r").Value2 = $nr.J
  $wsWriting.Range("K$r").Value2 = $nr.K
  $wsWriting.Range("L$r").Value2 = $nr.L

  $wsWriting.Range("M$r").Formula = "=SUM(C$r`:L$r)"
  $wsWriting.Range("N$r").FormulaArray = "=SUM(ABS(C$r`:L$r-C$prevRow`:L$prevRow))"
}

# Grow the table (and its AutoFilter) to cover the new rows.
$tbl = $wsWriting.ListObjects.Item(1)
$tbl.Resize($wsWriting.Range("A1:N42"))

# Recalculate so M/N show real totals instead of 0s.
$excel.Calculate()

# Extend the two chart series (bar = Daily/N, line = Total/M) to the new range.
$chart = $wsDashboard.ChartObjects(1).Chart
$chart.SeriesCollection(1).Formula = "=SERIES(writing!`$N`$1,writing!`$A`$2:`$A`$42,writing!`$N`$2:`$N`$42,1)"
$chart.SeriesCollection(2).Formula = "=SERIES(writing!`$M`$1,writing!`$A`$2:`$A`$42,writing!`$M`$2:`$M`$42,2)"

# Resize the chart (its anchor grows from col17/row27 to col27/row30).
$co = $wsDashboard.ChartObjects(1)
$co.Width  = 1577.8125
$co.Height = 427.5

# Restore the on-screen selections: "writing" remembers M41 (not the
# active tab), "dashboard" is re-selected last so it stays the active tab
# with its own new selection P33.
$wsWriting.Range("M41").Select()
$wsWriting.Application.ActiveWindow.ScrollRow = 17

$wsDashboard.Activate()
$wsDashboard.Range("P33").Select()
